# Append the 2025-03-19 price row (row 18) to every price sheet in the
# Solar_Prices workbook, mirroring the prior day's "Price" text formatting.
#
# Each worksheet holds a two-column Date/Price table where both columns are
# stored as text (inline strings) even though the Price column often looks
# numeric. We re-create that by entering the values with a leading
# apostrophe (forces Excel to store them as text) and then resetting the
# cell style to "Normal" so no stray quote-prefix/number-format style is
# left behind on the new cells.

$wb = $excel.ActiveWorkbook

$newDate = "2025-03-19"
$row = 18

$sheetValues = @{
    "N-Dense"                   = "40"
    "N-Type"                    = "43"
    "N-type Wafer"               = "1.19"
    "Cell Topcon 183mm"          = "0.298"
    "Module Topcon 183mm"        = "0.1"
    "Silver Rear_side"           = "5,450"
    "Silver Busbar front-side"   = "8,160"
    "Silver finger front-side"   = "8,210"
    "USD_CNY"                    = "7.2446"
}

foreach ($sheet in $wb.Worksheets) {
    $price = $sheetValues[$sheet.Name]
    if ($null -eq $price) {
        continue
    }

    $dateCell = $sheet.Cells.Item($row, 1)
    $priceCell = $sheet.Cells.Item($row, 2)

    $dateCell.Value = "'" + $newDate
    $priceCell.Value = "'" + $price

    $dateCell.Style = "Normal"
    $priceCell.Style = "Normal"
}
